$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Shift existing data rows 2-9 down to 3-10 one cell at a time (bottom-up,
# so a row is never overwritten before it has been read). `.Formula` is
# used to read the source value back reliably (plain `.Value` reads are
# not usable in this host), and the destination keeps/receives the same
# number format as the date column (D) relies on for its display.
for ($r = 9; $r -ge 2; $r--) {
    foreach ($col in $cols) {
        $ws.Range($col + ($r + 1)).Value = $ws.Range($col + $r).Formula
    }
}

# The date column (D) carries a custom date number format; make sure the
# freshly-written D10 (previously outside the used range) has it too.
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate row 2 with this week's new data point.
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C2").Value = "Ñuble"
$ws.Range("D2").Value = 45092
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100104
$ws.Range("H2").Value = "Frutos de pepita"
$ws.Range("I2").Value = 100104001
$ws.Range("J2").Value = "Granada"
$ws.Range("K2").Value = "Wonderfull"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 19000
$ws.Range("P2").Value = 18667
$ws.Range("Q2").Value = "$/caja 18 kilos granel"
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 1037
$ws.Range("T2").Value = 18
